# Add a new "Door V2" translation row to the Translations sheet.
#
# The new row is inserted right after the existing "object.LASER_FIELD"
# row (row 42), pushing every following row down by one. This also
# appends three new shared strings: "object.DOOR_V2", "Door V2" and
# "Puerta V2".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 42 (shifts rows 42..107 down to 43..108).
$ws.Range("A42").EntireRow.Insert()

# Populate the new row with the Door V2 translation strings.
$ws.Range("A42").Value = "object.DOOR_V2"
$ws.Range("B42").Value = "Door V2"
$ws.Range("C42").Value = "Puerta V2"

# Copy the formatting (style) used by the surrounding rows so the new
# row matches the rest of the table (left/top aligned, wrap text).
$ws.Range("A41:C41").Copy()
$ws.Range("A42:C42").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Update the view state to match the edited worksheet: scroll so row 23
# is at the top, and leave the selection on D42 (matching the saved
# view after editing row 42).
$excel.Goto($ws.Range("A23"), $true)
$ws.Range("D42").Select()
